$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel
# auto-converting numeric-looking strings (e.g. "212.27") into numbers.
# Builds the text via a text-formula in a scratch cell, then pastes the
# *value* of that formula (a String) into the target - PasteSpecial of a
# formula-produced string preserves the String type instead of re-parsing it.
function Set-TextValue($range, [string]$text) {
    $ws.Range("Z1").Formula = '="' + $text + '"'
    $ws.Range("Z1").Copy() | Out-Null
    $ws.Range($range).PasteSpecial(-4163) | Out-Null
}

Set-TextValue "D2" "27.616.67"
$ws.Range("E2").Value = "  +0.00%  "
Set-TextValue "D3" "1.631.03"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue "D5" "212.27"
$ws.Range("E5").Value = "  -0.13%  "
Set-TextValue "D6" "0.520"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  -0.13%  "
Set-TextValue "D8" "23.38"
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("E9").Value = "  +2.49%  "
$ws.Range("E10").Value = "  +0.20%  "
Set-TextValue "D11" "0.0874"
$ws.Range("E11").Value = "  -2.19%  "
Set-TextValue "D12" "1.860.69"
$ws.Range("E12").Value = "  -0.42%  "
Set-TextValue "D13" "1.636.83"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("E14").Value = "  +0.34%  "
Set-TextValue "D15" "0.553"
$ws.Range("E15").Value = "  -1.29%  "
Set-TextValue "D16" "65.37"
$ws.Range("E16").Value = "  +1.28%  "
Set-TextValue "D17" "27.584.89"
$ws.Range("E17").Value = "  -0.07%  "
Set-TextValue "D18" "231.60"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("E19").Value = "  -0.45%  "
Set-TextValue "D20" "7.58"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("E21").Value = "  -0.08%  "
Set-TextValue "D22" "10.57"
$ws.Range("E22").Value = "  +5.39%  "
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("E24").Value = "  +8.54%  "
Set-TextValue "D25" "149.58"
$ws.Range("E25").Value = "  -0.84%  "
Set-TextValue "D26" "6.90"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("E27").Value = "  -0.26%  "
Set-TextValue "D28" "15.55"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  -0.35%  "
Set-TextValue "D33" "1.475.33"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("E34").Value = "  -1.89%  "
Set-TextValue "D35" "1.56"
$ws.Range("E35").Value = "  -1.16%  "
Set-TextValue "D36" "2.34"
$ws.Range("E36").Value = "  -1.48%  "
Set-TextValue "D37" "0.940"
$ws.Range("E37").Value = "  +5.74%  "
Set-TextValue "D38" "0.880"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("E39").Value = "  +0.37%  "
Set-TextValue "D40" "0.557"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("E42").Value = "  -0.11%  "
Set-TextValue "D43" "68.01"
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("E44").Value = "  +0.45%  "
Set-TextValue "D45" "2.21"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("E46").Value = "  -4.37%  "
Set-TextValue "D47" "1.770.84"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  +1.66%  "
Set-TextValue "D49" "87.80"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").Value = "  +0.93%  "

# Clean up the scratch cell so it does not appear in the saved sheet.
$ws.Range("Z1").ClearContents()
